$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "tQfkE733"
$ws.Range("B2").Value = 23091138
$ws.Range("C2").Value = "dpjlgeo41"
$ws.Range("D2").Value = 'cp!5$H7A'
$ws.Range("F2").Value = "PWrcmZOI"
$ws.Range("G2").Value = "DMAu"
